$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Breast_Cancer_313)
$ws.Range("E2").Value = 234
$ws.Range("F2").Value = 74.76
$ws.Range("G2").Value = 172
$ws.Range("H2").Value = 54.95

# Row 3 (AABCGS_black)
$ws.Range("E3").Value = 87
$ws.Range("F3").Value = 97.75
$ws.Range("G3").Value = 37
$ws.Range("H3").Value = 41.57

# Row 4 (Gao2022_BrCa)
$ws.Range("E4").Value = 45386
$ws.Range("F4").Value = 79.7
$ws.Range("G4").Value = 41083
$ws.Range("H4").Value = 72.15

# Row 5 (Gao2022_ERNEG)
$ws.Range("E5").Value = 22608
$ws.Range("F5").Value = 79.59
$ws.Range("G5").Value = 20415
$ws.Range("H5").Value = 71.87

# Row 6 (Gao2022_ERPOS)
$ws.Range("E6").Value = 23370
$ws.Range("F6").Value = 79.74
$ws.Range("G6").Value = 21177
$ws.Range("H6").Value = 72.25

# Row 7 (Shieh2023)
$ws.Range("E7").Value = 66
$ws.Range("F7").Value = 86.84
